$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# like "332.70" or "1.002" are not silently coerced into floats
# (which would drop significant trailing zeros / add FP noise).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.873.19'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '1.747.32'
$ws.Range('E3').Value = '  -2.04%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '332.70'
$ws.Range('E5').Value = '  -1.51%  '
$ws.Range('D6').Value = '0.9952'
$ws.Range('E6').Value = '  -0.45%  '
$ws.Range('D7').Value = '0.3856'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').Value = '0.3378'
$ws.Range('E8').Value = '  -1.99%  '
$ws.Range('D9').Value = '45.38'
$ws.Range('E9').Value = '  -3.69%  '
$ws.Range('D10').Value = '1.106'
$ws.Range('E10').Value = '  -4.54%  '
$ws.Range('D11').Value = '0.07166'
$ws.Range('E11').Value = '  -3.65%  '
$ws.Range('D12').Value = '0.9996'
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '22.06'
$ws.Range('E13').Value = '  -6.20%  '
$ws.Range('D14').Value = '6.110'
$ws.Range('E14').Value = '  -5.15%  '
$ws.Range('D15').Value = '1.744.74'
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = '6.980'
$ws.Range('E16').Value = '  -3.36%  '
$ws.Range('D17').Value = '0.00001051'
$ws.Range('E17').Value = '  -2.44%  '
$ws.Range('D18').Value = '0.06586'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('D19').Value = '79.88'
$ws.Range('E19').Value = '  -3.36%  '
$ws.Range('D20').Value = '0.9961'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '16.80'
$ws.Range('E21').Value = '  -4.26%  '
$ws.Range('D22').Value = '6.177'
$ws.Range('E22').Value = '  -4.46%  '
$ws.Range('D23').Value = '27.919.32'
$ws.Range('E23').Value = '  -1.11%  '
$ws.Range('D24').Value = '11.51'
$ws.Range('E24').Value = '  -5.06%  '
$ws.Range('D25').Value = '2.377'
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('D26').Value = '153.73'
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('D27').Value = '19.80'
$ws.Range('E27').Value = '  -5.20%  '
$ws.Range('D28').Value = '2.295'
$ws.Range('E28').Value = '  -5.48%  '
$ws.Range('D29').Value = '1.950.20'
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').Value = '1.262'
$ws.Range('E30').Value = '  -12.81%  '
$ws.Range('D31').Value = '128.05'
$ws.Range('E31').Value = '  -7.05%  '
$ws.Range('D32').Value = '4.079'
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('D33').Value = '5.776'
$ws.Range('E33').Value = '  -6.26%  '
$ws.Range('D34').Value = '0.08663'
$ws.Range('E34').Value = '  -2.92%  '
$ws.Range('D35').Value = '12.00'
$ws.Range('E35').Value = '  -6.58%  '
$ws.Range('D36').Value = '0.02270'
$ws.Range('E36').Value = '  -6.79%  '
$ws.Range('D37').Value = '5.100'
$ws.Range('E37').Value = '  -4.59%  '
$ws.Range('D38').Value = '0.06093'
$ws.Range('E38').Value = '  -4.36%  '
$ws.Range('D39').Value = '1.504'
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('D40').Value = '0.6426'
$ws.Range('E40').Value = '  -6.49%  '
$ws.Range('D41').Value = '0.2091'
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').Value = '1.194'
$ws.Range('E42').Value = '  -3.77%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '7.911'
$ws.Range('E43').Value = '  -4.87%  '
$ws.Range('B44').Value = 'Frax'
$ws.Range('C44').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D44').Value = '0.9965'
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('D45').Value = '13.72'
$ws.Range('E45').Value = '  -3.94%  '
$ws.Range('D46').Value = '3.801'
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('D47').Value = '0.5939'
$ws.Range('E47').Value = '  -5.86%  '
$ws.Range('D48').Value = '125.80'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('D49').Value = '1.971'
$ws.Range('E49').Value = '  -6.04%  '
$ws.Range('D50').Value = '0.06960'
$ws.Range('E50').Value = '  -7.03%  '
$ws.Range('D51').Value = '1.144'
$ws.Range('E51').Value = '  -5.95%  '

# Restore the column back to its original (default/General) style so
# the only thing that changed in the saved file is cell content.
$ws.Range("D2:D51").Style = "Normal"

